$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns F and G
$ws.Range("F1").Value = "Hombre"
$ws.Range("G1").Value = "Mujer"

# Formulas for row 2 (Intercept)
$ws.Range("F2").Formula = "=B2"
$ws.Range("G2").Formula = "=B2"

# Formulas for row 3 (SEXO)
$ws.Range("F3").Formula = "=B3*1"
$ws.Range("G3").Formula = "=B3*0"

# Formulas for row 4 (EDAD)
$ws.Range("F4").Formula = "=B4*50.27"
$ws.Range("G4").Formula = "=B4*50.27"

# Formulas for row 5 (NIVELESTENTREV)
$ws.Range("F5").Formula = "=B5*8.09"
$ws.Range("G5").Formula = "=B5*8.09"

# Formulas for row 6 (SEXO:EDAD)
$ws.Range("F6").Formula = "=B6*1*50.27"
$ws.Range("G6").Formula = "=B6*0*50.27"

# Formulas for row 7 (SEXO:NIVELESTENTREV)
$ws.Range("F7").Formula = "=B7*1*8.09"
$ws.Range("G7").Formula = "=B7*0*8.09"

# Sum row 8
$ws.Range("F8").Formula = "=SUM(F2:F7)"
$ws.Range("G8").Formula = "=SUM(G2:G7)"

# Apply the same style as the rest of the header row (bold+centered) to F1:G1
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").HorizontalAlignment = -4108

# Match final selection/cursor position left by the author
$ws.Range("H8").Select() | Out-Null
